$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A88").Value = 229
$ws.Range("B88").Value = 138
$ws.Range("C88").Value = 84
$ws.Range("D88").Value = 4
$ws.Range("E88").Value = 3
$ws.Range("F88").Value = 87
$ws.Range("G88").Value = 88
$ws.Range("H88").Value = 8
$ws.Range("I88").Value = 0
